# ADDED REFERENCES FOR COMPONENTS ON SCHEMATIC
# Populate the "Schematic Reference Designators" column (Q) for each BOM row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$refs = @{
    3  = "D1"
    4  = "C9,C8,C7,C3"
    5  = "C1"
    6  = "C2"
    7  = "SW1, SW2"
    8  = "C10,C11"
    9  = "C6"
    10 = "L1"
    11 = "F1"
    12 = "D2"
    13 = "J5"
    14 = "J4"
    15 = "U3"
    16 = "U2"
    17 = "R4,R1"
    18 = "R3"
    19 = "J7,J6"
    20 = "D3"
    21 = "R2"
    22 = "D4,D5"
    23 = "U4,U5,U6,U7"
    24 = "J10,J8,J3"
}

foreach ($row in $refs.Keys) {
    $ws.Cells.Item($row, 17).Value = $refs[$row]
}
